$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = '68.480.74'
$r.Style = "Normal"
$ws.Range("E2").Value = '  -1.33%  '

$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = '2.451.18'
$r.Style = "Normal"
$ws.Range("E3").Value = '  -1.31%  '

$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '556.07'
$r.Style = "Normal"
$ws.Range("E5").Value = '  -1.99%  '

$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '161.69'
$r.Style = "Normal"
$ws.Range("E6").Value = '  -1.32%  '

$ws.Range("E7").Value = '  -0.05%  '

$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = '0.501'
$r.Style = "Normal"
$ws.Range("E8").Value = '  -1.95%  '

$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = '2.450.58'
$r.Style = "Normal"
$ws.Range("E9").Value = '  -1.35%  '

$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '0.148'
$r.Style = "Normal"
$ws.Range("E10").Value = '  -6.27%  '

$ws.Range("E11").Value = '  -1.23%  '

$ws.Range("E12").Value = '  -5.61%  '

$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = '4.77'
$r.Style = "Normal"
$ws.Range("E13").Value = '  -2.38%  '

$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '2.901.26'
$r.Style = "Normal"
$ws.Range("E14").Value = '  -1.41%  '

$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '68.375.30'
$r.Style = "Normal"
$ws.Range("E15").Value = '  -1.39%  '

$ws.Range("E16").Value = '  -4.09%  '

$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = '23.26'
$r.Style = "Normal"
$ws.Range("E17").Value = '  -3.63%  '

$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = '2.440.19'
$r.Style = "Normal"
$ws.Range("E18").Value = '  -1.90%  '

$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = '10.75'
$r.Style = "Normal"
$ws.Range("E19").Value = '  -3.27%  '

$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '339.63'
$r.Style = "Normal"
$ws.Range("E20").Value = '  -1.87%  '

$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '6.98'
$r.Style = "Normal"
$ws.Range("E21").Value = '  -5.07%  '

$ws.Range("E22").Value = '  -2.64%  '

$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("E24").Value = '  -2.49%  '

$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = '66.28'
$r.Style = "Normal"
$ws.Range("E25").Value = '  -4.71%  '

$ws.Range("E26").Value = '  -6.11%  '

$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = '2.577.91'
$r.Style = "Normal"
$ws.Range("E27").Value = '  -1.66%  '

$ws.Range("E28").Value = '  +0.19%  '

$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = '8.04'
$r.Style = "Normal"
$ws.Range("E29").Value = '  -6.36%  '

$ws.Range("E30").Value = '  -6.11%  '

$ws.Range("E31").Value = '  -6.31%  '

$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = '1.00'
$r.Style = "Normal"
$ws.Range("E32").Value = '  -0.03%  '

$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = '431.23'
$r.Style = "Normal"
$ws.Range("E33").Value = '  -0.90%  '

$ws.Range("E34").Value = '  -5.22%  '

$ws.Range("E35").Value = '  -5.49%  '

$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = '156.65'
$r.Style = "Normal"
$ws.Range("E36").Value = '  -0.10%  '

$ws.Range("E37").Value = '  -0.29%  '

$ws.Range("E38").Value = '  +0.12%  '

$ws.Range("E39").Value = '  -1.93%  '

$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '17.76'
$r.Style = "Normal"
$ws.Range("E40").Value = '  -1.83%  '

$ws.Range("E41").Value = '  -3.69%  '

$ws.Range("E42").Value = '  -3.24%  '

$ws.Range("E43").Value = '  -1.13%  '

$ws.Range("E44").Value = '  -7.16%  '

$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '1.08'
$r.Style = "Normal"
$ws.Range("E45").Value = '  +2.08%  '

$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '2.03'
$r.Style = "Normal"
$ws.Range("E46").Value = '  -5.34%  '

$ws.Range("E47").Value = '  -4.45%  '

$ws.Range("E48").Value = '  -2.67%  '

$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '0.0713'
$r.Style = "Normal"
$ws.Range("E49").Value = '  -1.30%  '

$ws.Range("E50").Value = '  -6.03%  '

$ws.Range("E51").Value = '  -2.11%  '
